# Updates the division-fact worksheet table: each cell's "A÷B=" expression
# is replaced with a new expression, per the target revision.
#
# Note on ordering: "68÷8=" is changed to "95÷8=", and separately the
# pre-existing "95÷8=" cell is changed to "78÷7=". To avoid the first
# replacement's output being re-matched by the second (or the second
# missing its original target because the first already overwrote it),
# the "95÷8=" -> "78÷7=" replacement runs BEFORE "68÷8=" -> "95÷8=".
# All other replacements are independent (no new value collides with any
# other old value) so their relative order does not matter.

$d = $word.ActiveDocument

$d.Content.Find.Execute("53÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷3=", 2) | Out-Null
$d.Content.Find.Execute("40÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷4=", 2) | Out-Null
$d.Content.Find.Execute("72÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=", 2) | Out-Null
$d.Content.Find.Execute("35÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷3=", 2) | Out-Null
$d.Content.Find.Execute("11÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷8=", 2) | Out-Null
$d.Content.Find.Execute("53÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷4=", 2) | Out-Null
$d.Content.Find.Execute("95÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷2=", 2) | Out-Null

# Must precede the "68÷8=" -> "95÷8=" replacement below (see note above).
$d.Content.Find.Execute("95÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷7=", 2) | Out-Null

$d.Content.Find.Execute("68÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=", 2) | Out-Null
$d.Content.Find.Execute("51÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷5=", 2) | Out-Null
$d.Content.Find.Execute("34÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=", 2) | Out-Null
$d.Content.Find.Execute("96÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷4=", 2) | Out-Null
$d.Content.Find.Execute("49÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷3=", 2) | Out-Null
$d.Content.Find.Execute("74÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷4=", 2) | Out-Null
$d.Content.Find.Execute("11÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷2=", 2) | Out-Null
$d.Content.Find.Execute("46÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷3=", 2) | Out-Null
$d.Content.Find.Execute("66÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷8=", 2) | Out-Null
$d.Content.Find.Execute("40÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷7=", 2) | Out-Null
$d.Content.Find.Execute("27÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷2=", 2) | Out-Null
$d.Content.Find.Execute("87÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷2=", 2) | Out-Null
$d.Content.Find.Execute("25÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷8=", 2) | Out-Null
$d.Content.Find.Execute("98÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷5=", 2) | Out-Null
$d.Content.Find.Execute("67÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷3=", 2) | Out-Null
$d.Content.Find.Execute("99÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷8=", 2) | Out-Null
$d.Content.Find.Execute("41÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷5=", 2) | Out-Null
